# Data update from DGS's 2021/08/11 report.
# Append a new row (66) to Sheet1 with the latest time-series entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 66

# Column A holds dates stored as literal text (matching every other row in
# this column), not real date serials. Writing the string directly makes
# Excel auto-detect it as a date and reformat the cell, so stage the text
# in a scratch cell that is pre-formatted as Text, then copy only the
# value (and its "is text" flag) into place - this preserves the existing
# date-style formatting (style index) already applied to column A.
$scratch = $ws.Cells.Item(200, 50)
$scratch.NumberFormat = "@"
$scratch.Value = "2021/08/11"
$scratch.Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()

$ws.Cells.Item($row, 2).Value = 326.5
$ws.Cells.Item($row, 3).Value = 331.6
$ws.Cells.Item($row, 4).Value = 0.94
$ws.Cells.Item($row, 5).Value = 0.94

# Move the active selection to the next empty row, like a user would see
# after entering a new line of data at the bottom of the table.
$null = $ws.Range("A67").Select()
